$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '62.832.22'
Set-TextValue $ws.Range("E2") '  +0.16%  '

Set-TextValue $ws.Range("D3") '2.461.32'
Set-TextValue $ws.Range("E3") '  +0.62%  '

Set-TextValue $ws.Range("E4") '  +0.06%  '

Set-TextValue $ws.Range("D5") '574.11'
Set-TextValue $ws.Range("E5") '  -0.35%  '

Set-TextValue $ws.Range("D6") '146.64'
Set-TextValue $ws.Range("E6") '  +0.66%  '

Set-TextValue $ws.Range("E7") '  +0.00%  '

Set-TextValue $ws.Range("D8") '0.535'
Set-TextValue $ws.Range("E8") '  -0.64%  '

Set-TextValue $ws.Range("D9") '2.461.73'
Set-TextValue $ws.Range("E9") '  +0.70%  '

Set-TextValue $ws.Range("E10") '  +0.36%  '

Set-TextValue $ws.Range("E11") '  -0.35%  '

Set-TextValue $ws.Range("E12") '  +0.55%  '

Set-TextValue $ws.Range("D13") '0.356'
Set-TextValue $ws.Range("E13") '  +0.99%  '

Set-TextValue $ws.Range("D14") '29.04'
Set-TextValue $ws.Range("E14") '  +3.18%  '

Set-TextValue $ws.Range("E15") '  -0.61%  '

Set-TextValue $ws.Range("D16") '2.909.11'
Set-TextValue $ws.Range("E16") '  +0.70%  '

Set-TextValue $ws.Range("D17") '62.782.97'
Set-TextValue $ws.Range("E17") '  +0.14%  '

Set-TextValue $ws.Range("D18") '2.467.56'
Set-TextValue $ws.Range("E18") '  +0.71%  '

Set-TextValue $ws.Range("D19") '7.93'
Set-TextValue $ws.Range("E19") '  +0.43%  '

Set-TextValue $ws.Range("D20") '10.98'
Set-TextValue $ws.Range("E20") '  -0.01%  '

Set-TextValue $ws.Range("D21") '326.34'
Set-TextValue $ws.Range("E21") '  -0.91%  '

Set-TextValue $ws.Range("E22") '  -0.13%  '

Set-TextValue $ws.Range("D23") '2.21'
Set-TextValue $ws.Range("E23") '  +8.08%  '

Set-TextValue $ws.Range("E24") '  -0.01%  '

Set-TextValue $ws.Range("D25") '9.95'
Set-TextValue $ws.Range("E25") '  +16.57%  '

Set-TextValue $ws.Range("D26") '65.55'
Set-TextValue $ws.Range("E26") '  -0.85%  '

Set-TextValue $ws.Range("D27") '647.14'
Set-TextValue $ws.Range("E27") '  -0.22%  '

Set-TextValue $ws.Range("D28") '0.0₃0985'
Set-TextValue $ws.Range("E28") '  -0.53%  '

Set-TextValue $ws.Range("D29") '2.591.66'

Set-TextValue $ws.Range("D30") '0.996'
Set-TextValue $ws.Range("E30") '  -15.34%  '

Set-TextValue $ws.Range("E31") '  -1.20%  '

Set-TextValue $ws.Range("D32") '7.97'
Set-TextValue $ws.Range("E32") '  -2.81%  '

Set-TextValue $ws.Range("E33") '  -1.53%  '

Set-TextValue $ws.Range("E34") '  -2.98%  '

Set-TextValue $ws.Range("E35") '  -0.03%  '

Set-TextValue $ws.Range("E36") '  +2.86%  '

Set-TextValue $ws.Range("E37") '  -0.46%  '

Set-TextValue $ws.Range("D38") '2.84'
Set-TextValue $ws.Range("E38") '  +4.27%  '

Set-TextValue $ws.Range("D39") '0.368'
Set-TextValue $ws.Range("E39") '  -1.53%  '

Set-TextValue $ws.Range("B40") 'Monero'
Set-TextValue $ws.Range("C40") 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D40") '151.48'
Set-TextValue $ws.Range("E40") '  -1.36%  '

Set-TextValue $ws.Range("B41") 'EthereumClassic'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D41") '18.69'
Set-TextValue $ws.Range("E41") '  -0.31%  '

Set-TextValue $ws.Range("B42") 'RenderToken'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range("D42") '5.37'
Set-TextValue $ws.Range("E42") '  -2.56%  '

Set-TextValue $ws.Range("E43") '  -1.22%  '

Set-TextValue $ws.Range("D44") '0.0₆0307'
Set-TextValue $ws.Range("E44") '  -46.77%  '

Set-TextValue $ws.Range("E45") '  +0.00%  '

Set-TextValue $ws.Range("D46") '152.30'
Set-TextValue $ws.Range("E46") '  +4.64%  '

Set-TextValue $ws.Range("E47") '  +2.21%  '

Set-TextValue $ws.Range("D48") '3.57'
Set-TextValue $ws.Range("E48") '  -1.71%  '

Set-TextValue $ws.Range("D49") '20.44'
Set-TextValue $ws.Range("E49") '  -1.34%  '

Set-TextValue $ws.Range("D50") '0.607'
Set-TextValue $ws.Range("E50") '  +0.26%  '

Set-TextValue $ws.Range("D51") '0.0510'
Set-TextValue $ws.Range("E51") '  -1.47%  '
